$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing column A/B data: B3 changes from 29 to 20
$ws.Range("B3").Value = 20

# B4 gets an underline font style but stays empty
$ws.Range("B4").Font.Underline = $true

# New column D/E block
$ws.Range("D1").Value = "Origem lead - 11"
$ws.Range("E1").HorizontalAlignment = -4108  # xlCenter

$ws.Range("D2").Value = "PRISCYLLA"
$ws.Range("E2").Value = 1
$ws.Range("E2").Font.Underline = $true
$ws.Range("E2").HorizontalAlignment = -4108

$ws.Range("D3").Value = "PROSP MAURO"
$ws.Range("E3").Value = 2
$ws.Range("E3").HorizontalAlignment = -4108

$ws.Range("D4").Value = "VISÃO CEDENTE"
$ws.Range("E4").Value = 8
$ws.Range("E4").HorizontalAlignment = -4108

$ws.Range("D5").Value = "ALEX"
$ws.Range("E5").Value = 0
$ws.Range("E5").HorizontalAlignment = -4108

# Selection / view adjustments
$ws.Range("G3").Select()
